# Weekly driver report update for 2025-04-20
#
# The "Good Drivers" table (rows 12-19 on the Driver Summary sheet) is
# re-sorted by Driver Vintage date (most recent first) and refreshed with
# this week's sample counts / vintage dates. Row 12 previously had no
# vintage date and now becomes the newest entry; row 13 gains a vintage
# date for the first time; row 15's client-count is bumped; the rest of
# the rows simply shift position to stay sorted by date.
#
# Column E ("Driver Vintage") stores plain text dates (e.g. "2021-01-19"),
# not real Excel date serials, so each value is entered with a leading
# apostrophe to force literal text and stop Excel's automatic date
# recognition from converting it into a date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Intel(R) Dual Band Wireless-AC 8265 - 22.30.0.11
$ws.Range("A12").Value = "Intel(R) Dual Band Wireless-AC 8265 - 22.30.0.11"
$ws.Range("B12").Value = 170510
$ws.Range("D12").Value = 99.90000000000001
$ws.Range("E12").Value = "'2021-01-19"

# Row 13: Intel(R) Dual Band Wireless-AC 8265 - 22.0.1.1
$ws.Range("A13").Value = "Intel(R) Dual Band Wireless-AC 8265 - 22.0.1.1"
$ws.Range("B13").Value = 52096
$ws.Range("D13").Value = 100
$ws.Range("E13").Value = "'2020-09-28"

# Row 14: Intel(R) Dual Band Wireless-AC 8265 - 20.70.11.3
$ws.Range("A14").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.11.3"
$ws.Range("B14").Value = 161874
$ws.Range("D14").Value = 100
$ws.Range("E14").Value = "'2019-09-05"

# Row 15: Intel(R) Dual Band Wireless-AC 8265 - 20.70.12.5
$ws.Range("A15").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.12.5"
$ws.Range("B15").Value = 143342
$ws.Range("D15").Value = 99.90000000000001
$ws.Range("E15").Value = "'2019-08-25"

# Row 16: Intel(R) Dual Band Wireless-AC 8265 - 20.70.10.2
$ws.Range("A16").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.10.2"
$ws.Range("B16").Value = 20227
$ws.Range("D16").Value = 100
$ws.Range("E16").Value = "'2019-05-11"

# Row 17: Intel(R) Dual Band Wireless-AC 8265 - 20.70.9.1
$ws.Range("A17").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.9.1"
$ws.Range("B17").Value = 34065
$ws.Range("D17").Value = 100
$ws.Range("E17").Value = "'2019-04-28"

# Row 18: Intel(R) Dual Band Wireless-AC 8265 - 20.70.8.1
$ws.Range("A18").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.8.1"
$ws.Range("B18").Value = 48540
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "'2019-03-16"

# Row 19: Intel(R) Dual Band Wireless-AC 8265 - 20.70.5.2
$ws.Range("A19").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.5.2"
$ws.Range("B19").Value = 184564
$ws.Range("D19").Value = 99.90000000000001
$ws.Range("E19").Value = "'2018-11-25"
